$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Controller" (B) and "Five Year Plan" (C) values for fiscal years 2017-2028
$ws.Range("B2").Value = 31942427.94974807
$ws.Range("C2").Value = 31942427.94974807

$ws.Range("B3").Value = 78152731.89953387
$ws.Range("C3").Value = 78152731.89953387

$ws.Range("B4").Value = 77392638.51289788
$ws.Range("C4").Value = 77392638.51289788

$ws.Range("B5").Value = 70685021.27613598
$ws.Range("C5").Value = 70685021.27613598

$ws.Range("B6").Value = 69527511.64367333
$ws.Range("C6").Value = 69527511.64367333

$ws.Range("B7").Value = 75093575.51524313
$ws.Range("C7").Value = 75093575.51524313

$ws.Range("B8").Value = 72594000
$ws.Range("C8").Value = 72594000

$ws.Range("B9").Value = 73579540.96680468
$ws.Range("C9").Value = 73501000

$ws.Range("B10").Value = 72680917.20352739
$ws.Range("C10").Value = 74368000

$ws.Range("B11").Value = 71921422.15560062
$ws.Range("C11").Value = 73624000

$ws.Range("B12").Value = 71161594.22186247
$ws.Range("C12").Value = 73226000

$ws.Range("B13").Value = 70520384.36223672
$ws.Range("C13").Value = 72911000
